$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New rows 168-174 continuing the microciclo schedule (Microciclo 25)
$data = @(
    @{ Row = 168; Date = "2025-11-03"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "DESCANSO"; Int = $null; Partido = $null },
    @{ Row = 169; Date = "2025-11-04"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "ENTRENO";   Int = 1;     Partido = $null },
    @{ Row = 170; Date = "2025-11-05"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "ENTRENO";   Int = 2;     Partido = $null },
    @{ Row = 171; Date = "2025-11-06"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "ENTRENO";   Int = -2;    Partido = $null },
    @{ Row = 172; Date = "2025-11-07"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "ENTRENO";   Int = -1;    Partido = $null },
    @{ Row = 173; Date = "2025-11-08"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "PARTIDO";   Int = $null; Partido = "Monterrey" },
    @{ Row = 174; Date = "2025-11-09"; Micro = 25; Tipo = "Competencia"; Fase = "Competencia"; Dia = "DESCANSO"; Int = $null; Partido = $null }
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 1).Value = [DateTime]$r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Micro
    $ws.Cells.Item($r.Row, 3).Value = $r.Tipo
    $ws.Cells.Item($r.Row, 4).Value = $r.Fase
    $ws.Cells.Item($r.Row, 5).Value = $r.Dia
    if ($null -ne $r.Int) {
        $ws.Cells.Item($r.Row, 6).Value = $r.Int
    }
    if ($null -ne $r.Partido) {
        $ws.Cells.Item($r.Row, 7).Value = $r.Partido
    }
}

$ws.Range("A175:D175").Select()
